# Update crypto price/volume figures to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the cell to hold the exact literal string (no numeric
    # auto-coercion of values like "593.96"), then drop back to the
    # workbook default style so no stray number-format/style sticks.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '62.847.30'
Set-TextCell $ws.Range("E2") '  +2.12%  '

Set-TextCell $ws.Range("D3") '3.034.24'
Set-TextCell $ws.Range("E3") '  +1.29%  '

Set-TextCell $ws.Range("D5") '593.96'
Set-TextCell $ws.Range("E5") '  -0.36%  '

Set-TextCell $ws.Range("D6") '153.61'
Set-TextCell $ws.Range("E6") '  +6.81%  '

Set-TextCell $ws.Range("D8") '3.028.56'
Set-TextCell $ws.Range("E8") '  +1.13%  '

Set-TextCell $ws.Range("E9") '  -0.61%  '

Set-TextCell $ws.Range("D10") '6.53'
Set-TextCell $ws.Range("E10") '  +10.53%  '

Set-TextCell $ws.Range("D11") '0.151'
Set-TextCell $ws.Range("E11") '  +2.40%  '

Set-TextCell $ws.Range("D12") '0.465'
Set-TextCell $ws.Range("E12") '  +1.03%  '

Set-TextCell $ws.Range("D13") '0.0000234'
Set-TextCell $ws.Range("E13") '  +2.36%  '

Set-TextCell $ws.Range("D14") '35.65'
Set-TextCell $ws.Range("E14") '  +3.74%  '

Set-TextCell $ws.Range("E15") '  +2.12%  '

Set-TextCell $ws.Range("D16") '3.539.96'
Set-TextCell $ws.Range("E16") '  +1.40%  '

Set-TextCell $ws.Range("D17") '7.13'
Set-TextCell $ws.Range("E17") '  +1.40%  '

Set-TextCell $ws.Range("D18") '62.892.06'
Set-TextCell $ws.Range("E18") '  +2.25%  '

Set-TextCell $ws.Range("D19") '3.035.02'
Set-TextCell $ws.Range("E19") '  +1.25%  '

Set-TextCell $ws.Range("D20") '453.08'
Set-TextCell $ws.Range("E20") '  -0.05%  '

Set-TextCell $ws.Range("D21") '14.29'
Set-TextCell $ws.Range("E21") '  +2.26%  '

Set-TextCell $ws.Range("D22") '0.697'
Set-TextCell $ws.Range("E22") '  +1.53%  '

Set-TextCell $ws.Range("D23") '7.48'
Set-TextCell $ws.Range("E23") '  +1.77%  '

Set-TextCell $ws.Range("D24") '83.12'
Set-TextCell $ws.Range("E24") '  +1.45%  '

Set-TextCell $ws.Range("D25") '11.33'
Set-TextCell $ws.Range("E25") '  +7.48%  '

Set-TextCell $ws.Range("D26") '2.31'
Set-TextCell $ws.Range("E26") '  +3.56%  '

Set-TextCell $ws.Range("E27") '  +3.38%  '

Set-TextCell $ws.Range("E28") '  +0.06%  '

Set-TextCell $ws.Range("D29") '7.42'
Set-TextCell $ws.Range("E29") '  +3.02%  '

Set-TextCell $ws.Range("E30") '  +0.99%  '

Set-TextCell $ws.Range("E31") '  +6.59%  '

Set-TextCell $ws.Range("E32") '  -0.04%  '

Set-TextCell $ws.Range("D33") '27.58'
Set-TextCell $ws.Range("E33") '  +0.37%  '

Set-TextCell $ws.Range("E34") '  +1.71%  '

Set-TextCell $ws.Range("D35") '0.0₃0866'
Set-TextCell $ws.Range("E35") '  +4.17%  '

Set-TextCell $ws.Range("E36") '  +2.02%  '

Set-TextCell $ws.Range("E37") '  +2.92%  '

Set-TextCell $ws.Range("D38") '3.21'
Set-TextCell $ws.Range("E38") '  +11.40%  '

Set-TextCell $ws.Range("D39") '2.11'
Set-TextCell $ws.Range("E39") '  +2.41%  '

Set-TextCell $ws.Range("D40") '50.56'
Set-TextCell $ws.Range("E40") '  +0.38%  '

Set-TextCell $ws.Range("E41") '  +3.51%  '

Set-TextCell $ws.Range("D42") '9.11'
Set-TextCell $ws.Range("E42") '  -1.95%  '

Set-TextCell $ws.Range("D43") '0.304'
Set-TextCell $ws.Range("E43") '  +12.87%  '

Set-TextCell $ws.Range("D44") '41.84'
Set-TextCell $ws.Range("E44") '  +6.42%  '

Set-TextCell $ws.Range("D45") '398.28'
Set-TextCell $ws.Range("E45") '  +0.66%  '

Set-TextCell $ws.Range("E46") '  +1.91%  '

Set-TextCell $ws.Range("D47") '2.730.07'
Set-TextCell $ws.Range("E47") '  +0.50%  '

Set-TextCell $ws.Range("D48") '132.07'
Set-TextCell $ws.Range("E48") '  -0.91%  '

Set-TextCell $ws.Range("E49") '  +0.05%  '

Set-TextCell $ws.Range("D50") '2.27'
Set-TextCell $ws.Range("E50") '  +5.18%  '

Set-TextCell $ws.Range("D51") '24.47'
Set-TextCell $ws.Range("E51") '  +4.43%  '
